$wb = $excel.ActiveWorkbook

$wsWorkers   = $wb.Worksheets.Item("Workers")
$wsBios      = $wb.Worksheets.Item("Bios")
$wsSkills    = $wb.Worksheets.Item("Skills")
$wsContracts = $wb.Worksheets.Item("Contracts")

# ---------------------------------------------------------------------------
# Workers sheet (worker renamed Clayton -> Beebee, plus several attribute
# changes: UID, birthday/debut date, body type, height/weight range, race,
# style, gimmick text and bio).
# ---------------------------------------------------------------------------
$wsWorkers.Range("A2").Value  = 2939
$wsWorkers.Range("E2").Value  = "Beebee"
$wsWorkers.Range("F2").Value  = "Beebee"
$wsWorkers.Range("L2").Value  = 24401
$wsWorkers.Range("M2").Value  = 34560
$wsWorkers.Range("O2").Value  = 1
$wsWorkers.Range("P2").Value  = 34
$wsWorkers.Range("Q2").Value  = 270
$wsWorkers.Range("R2").Value  = 250
$wsWorkers.Range("S2").Value  = 309
$wsWorkers.Range("T2").Value  = "beebee.jpg"
$wsWorkers.Range("V2").Value  = 6
$wsWorkers.Range("AC2").Value = 11
$wsWorkers.Range("BL2").Value = 'Wrestling Gimmick: "The Sultan'
$wsWorkers.Range("BN2").Value = "Wrestling Gimmick: The Siniste"

# ---------------------------------------------------------------------------
# Bios sheet
# ---------------------------------------------------------------------------
$wsBios.Range("A2").Value = 2939

$newBio = @"
Beebee, hailing from the Middle East, is a professional wrestler who has been captivating audiences around the world with his unique style and charisma. With his mysterious aura and enigmatic presence, Beebee brings a sense of intrigue to the wrestling ring like no other.
Known for his Interpret style of wrestling, Beebee incorporates elements of storytelling and emotion into his matches, creating a truly captivating experience for fans. Whether he is portraying a hero fighting against insurmountable odds or a villain seeking to undermine his opponents, Beebee always leaves a lasting impression with his performances.
Despite his enigmatic persona, Beebee is a fierce competitor in the ring, utilizing his agility, strength, and technical prowess to outsmart and outmaneuver his opponents. With a reputation for delivering show-stopping matches and electrifying moments, Beebee has quickly risen through the ranks of professional wrestling and established himself as a force to be reckoned with.
Off the mat, Beebee is a dedicated athlete and performer who is constantly honing his craft and pushing the boundaries of what is possible in the world of professional wrestling. With his unwavering dedication to his art and his undeniable talent, Beebee is sure to continue making waves in the wrestling world for years to come.
"@
$wsBios.Range("B2").Value = $newBio.TrimEnd("`r", "`n")

# ---------------------------------------------------------------------------
# Skills sheet
# ---------------------------------------------------------------------------
$wsSkills.Range("A2").Value  = 2939
$wsSkills.Range("B2").Value  = 71
$wsSkills.Range("C2").Value  = 35
$wsSkills.Range("D2").Value  = 29
$wsSkills.Range("E2").Value  = 63
$wsSkills.Range("F2").Value  = 65
$wsSkills.Range("G2").Value  = 23
$wsSkills.Range("H2").Value  = 61
$wsSkills.Range("J2").Value  = 72
$wsSkills.Range("K2").Value  = 86
$wsSkills.Range("L2").Value  = 71
$wsSkills.Range("M2").Value  = 41
$wsSkills.Range("N2").Value  = 49
$wsSkills.Range("Q2").Value  = 49
$wsSkills.Range("R2").Value  = 24
$wsSkills.Range("S2").Value  = 61
$wsSkills.Range("T2").Value  = 76
$wsSkills.Range("U2").Value  = 53
$wsSkills.Range("V2").Value  = 48
$wsSkills.Range("W2").Value  = 69
$wsSkills.Range("X2").Value  = 22
$wsSkills.Range("Y2").Value  = 37

# ---------------------------------------------------------------------------
# Contracts sheet
# ---------------------------------------------------------------------------
$wsContracts.Range("A2").Value  = 2628
$wsContracts.Range("B2").Value  = 38
$wsContracts.Range("C2").Value  = 2939
$wsContracts.Range("H2").Value  = $true
$wsContracts.Range("AA2").Value = 41557
$wsContracts.Range("AB2").Value = 23
$wsContracts.Range("AJ2").Value = 5
$wsContracts.Range("AX2").Value = 1
